$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new column before column A. This shifts the existing
#    columns (query/StatQuery/dbExcel/WebExcel + their row-2 values)
#    one slot to the right (A->B, B->C, C->D, D->E) and preserves the
#    exact stored column widths and the wrap-text style on the old
#    A2 cell (which becomes B2).
# ------------------------------------------------------------------
$ws.Columns.Item(1).Insert()

# New column A width (TabName column).
$ws.Columns.Item(1).ColumnWidth = 10.0

# The column that used to be "B" (StatQuery) is now "C" - it was a
# bestFit 255-wide column; narrow it down now that its longest value
# is shorter.
$ws.Columns.Item(3).ColumnWidth = 91.8

# ------------------------------------------------------------------
# 2. New column A values - a "TabName" label column.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# ------------------------------------------------------------------
# 3. Row 2 (Cases tab) - replace the case query in column B with the
#    new, richer query text; column C keeps the StatQuery text
#    (already shifted there by the column insert) but now also needs
#    the wrap-text formatting that used to live only on column A.
# ------------------------------------------------------------------
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Basset Hound']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@
$ws.Range("B2").Value = $casesQuery
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

# ------------------------------------------------------------------
# 4. Row 3 (Samples tab) - brand new row.
# ------------------------------------------------------------------
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Basset Hound']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$statQuery = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Basset Hound']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = "TC04_Canine_Filter_Breed-BassHnd_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC04_Canine_Filter_Breed-BassHnd_WebData.xlsx"
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# ------------------------------------------------------------------
# 5. Row 4 (Files tab) - brand new row.
# ------------------------------------------------------------------
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Basset Hound']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = "TC04_Canine_Filter_Breed-BassHnd_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC04_Canine_Filter_Breed-BassHnd_WebData.xlsx"
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# ------------------------------------------------------------------
# 6. Row heights for the (now much longer) wrapped query cells.
# ------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 275.5
$ws.Rows.Item(3).RowHeight = 232
$ws.Rows.Item(4).RowHeight = 246.5

# ------------------------------------------------------------------
# 7. View state - zoomed out a bit and scrolled/selected differently.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 40
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("J4").Select()
